$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before D, shifting existing D:K data to E:L
# (2018 fiscal year column added ahead of the existing years)
$ws.Columns("D:D").Insert()

# Copy number formats/styles from the (now shifted) column E into the
# newly inserted column D so the new cells match the existing look
# (date format for the header row, number format for data rows, etc.)
$ws.Columns("E:E").Copy()
$ws.Columns("D:D").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the 2018 fiscal year figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 533300
$ws.Range("D9").Value = 119800
$ws.Range("D10").Value = 413400
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 49400
$ws.Range("D17").Value = 491700
$ws.Range("D18").Value = 41600
$ws.Range("D20").Value = -4000
$ws.Range("D21").Value = 87000
$ws.Range("D22").Value = 6200
$ws.Range("D23").Value = 31300
$ws.Range("D24").Value = 7300
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 24000
$ws.Range("D27").Value = 24000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 4000
$ws.Range("D33").Value = 24000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 24000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 300
$ws.Range("D42").Value = 27500
$ws.Range("D43").Value = 69000
$ws.Range("D44").Value = 1500
$ws.Range("D45").Value = 10400
$ws.Range("D46").Value = 108700
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 354600
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 2800
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 466100
$ws.Range("D57").Value = 20000
$ws.Range("D58").Value = 66600
$ws.Range("D59").Value = 20800
$ws.Range("D60").Value = 107400
$ws.Range("D61").Value = 157300
$ws.Range("D62").Value = 61900
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 326600
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 199100
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 139400
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 24000
$ws.Range("D83").Value = 49400
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 82300
$ws.Range("D91").Value = -73900
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -55300
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -27000
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 100

# Two small corrections to the existing (now shifted) 2017 column E
# values that were updated as part of this data refresh
$ws.Range("E89").Value = 50600
$ws.Range("E94").Value = -45300
